$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.95"
$ws.Range("D3").Value = "'21.73"
$ws.Range("D4").Value = "'5.429"
$ws.Range("D5").Value = "'0.05687"
$ws.Range("D6").Value = "'3.381"
$ws.Range("D7").Value = "'0.8080"
$ws.Range("D8").Value = "'1.030"
$ws.Range("D9").Value = "'0.1468"
$ws.Range("D10").Value = "'0.07784"
$ws.Range("D11").Value = "'0.03172"
$ws.Range("D12").Value = "'0.03064"
$ws.Range("D13").Value = "'0.09265"
$ws.Range("D14").Value = "'3.570"
$ws.Range("D15").Value = "'0.001655"
$ws.Range("D16").Value = "'0.04721"
$ws.Range("D17").Value = "'0.0005863"
$ws.Range("D18").Value = "'0.006360"
$ws.Range("D19").Value = "'0.005041"
$ws.Range("D20").Value = "'0.001044"
$ws.Range("D21").Value = "'0.0001502"
$ws.Range("D22").Value = "'0.0003201"
$ws.Range("D23").Value = "'3.771"
$ws.Range("D26").Value = "'0.3304"
$ws.Range("D40").Value = "'0.04075"
$ws.Range("D41").Value = "'0.006938"
$ws.Range("D42").Value = "'0.1046"
$ws.Range("D43").Value = "'0.003022"
$ws.Range("D44").Value = "'0.007750"
$ws.Range("D45").Value = "'0.00005905"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.0005503"
$ws.Range("D48").Value = "'0.6828"
$ws.Range("D49").Value = "'0.008891"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D51").Value = "'0.01010"
